$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, shifting existing rows 284-390 down to 285-391.
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row 284 with the new data record.
$ws.Range("A284").Value = 7
$ws.Range("B284").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C284").Value = "Ñuble"
$ws.Range("D284").Value = 45229
$ws.Range("E284").Value = 16
$ws.Range("F284").Value = 100112043
$ws.Range("G284").Value = "Pepino ensalada"
$ws.Range("H284").Value = "Sin especificar"
$ws.Range("I284").Value = "Primera"
$ws.Range("J284").Value = 100
$ws.Range("K284").Value = 14000
$ws.Range("L284").Value = 14000
$ws.Range("M284").Value = 14000
$ws.Range("N284").Value = "$/caja 60 unidades"
$ws.Range("O284").Value = "Región de Arica y Parinacota"
$ws.Range("P284").Value = 233
$ws.Range("Q284").Value = 60
$ws.Range("R284").Value = "Hortaliza"
